$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 128: swap in match data that was previously on row 129 (id in column A is unchanged)
$ws.Range("B128").Value = 7462738
$ws.Range("E128").Value = 'Vaca Diez'
$ws.Range("F128").Value = 'The Strongest'
$ws.Range("G128").Value = 2
$ws.Range("H128").Value = 2
$ws.Range("J128").Value = 1
$ws.Range("K128").Value = 'D'
$ws.Range("L128").Value = 4
$ws.Range("M128").Value = 4
$ws.Range("N128").Value = 1.666
$ws.Range("O128").Value = 4
$ws.Range("P128").Value = 3.8
$ws.Range("Q128").Value = 1.75
$ws.Range("R128").Value = 0.75
$ws.Range("S128").Value = 1.8
$ws.Range("T128").Value = 2
$ws.Range("U128").Value = 3
$ws.Range("V128").Value = 1.925
$ws.Range("W128").Value = 1.875
$ws.Range("X128").Value = -1
$ws.Range("Y128").Value = 2.8
$ws.Range("AA128").Value = 0.8
$ws.Range("AC128").Value = 0.925
$ws.Range("AD128").Value = -1

# Row 129: swap in match data that was previously on row 128 (id in column A is unchanged)
$ws.Range("B129").Value = 7462542
$ws.Range("E129").Value = 'Always Ready'
$ws.Range("F129").Value = 'Royal Pari FC'
$ws.Range("G129").Value = 3
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 'H'
$ws.Range("L129").Value = 1.363
$ws.Range("M129").Value = 4.75
$ws.Range("N129").Value = 6.5
$ws.Range("O129").Value = 1.285
$ws.Range("P129").Value = 6.5
$ws.Range("Q129").Value = 8
$ws.Range("R129").Value = -1.75
$ws.Range("S129").Value = 1.9
$ws.Range("T129").Value = 1.9
$ws.Range("U129").Value = 3.25
$ws.Range("V129").Value = 1.85
$ws.Range("W129").Value = 1.95
$ws.Range("X129").Value = 0.2849999999999999
$ws.Range("Y129").Value = -1
$ws.Range("AA129").Value = 0.8999999999999999
$ws.Range("AC129").Value = -0.5
$ws.Range("AD129").Value = 0.475

# Row 143: swap in match data that was previously on row 144 (id in column A is unchanged)
$ws.Range("B143").Value = 7532412
$ws.Range("E143").Value = 'Vaca Diez'
$ws.Range("F143").Value = 'Blooming'
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 3
$ws.Range("I143").Value = 0
$ws.Range("J143").Value = 2
$ws.Range("K143").Value = 'A'
$ws.Range("L143").Value = 1.727
$ws.Range("N143").Value = 4
$ws.Range("O143").Value = 2.3
$ws.Range("P143").Value = 3.6
$ws.Range("Q143").Value = 2.875
$ws.Range("R143").Value = -0.25
$ws.Range("S143").Value = 1.95
$ws.Range("T143").Value = 1.85
$ws.Range("U143").Value = 2.75
$ws.Range("V143").Value = 1.925
$ws.Range("W143").Value = 1.875
$ws.Range("X143").Value = -1
$ws.Range("Z143").Value = 1.875
$ws.Range("AB143").Value = 0.8500000000000001
$ws.Range("AC143").Value = 0.4625
$ws.Range("AD143").Value = -0.5

# Row 144: swap in match data that was previously on row 145 (id in column A is unchanged)
$ws.Range("B144").Value = 7532413
$ws.Range("E144").Value = 'Libertad Gran Mamore FC'
$ws.Range("F144").Value = 'Club Aurora'
$ws.Range("H144").Value = 1
$ws.Range("J144").Value = 0
$ws.Range("L144").Value = 2.25
$ws.Range("M144").Value = 3.3
$ws.Range("N144").Value = 2.8
$ws.Range("O144").Value = 2.375
$ws.Range("P144").Value = 3.4
$ws.Range("S144").Value = 2.025
$ws.Range("T144").Value = 1.775
$ws.Range("U144").Value = 2.5
$ws.Range("V144").Value = 1.9
$ws.Range("W144").Value = 1.9
$ws.Range("AB144").Value = 0.7749999999999999
$ws.Range("AC144").Value = -1
$ws.Range("AD144").Value = 0.8999999999999999

# Row 145: swap in match data that was previously on row 143 (id in column A is unchanged)
$ws.Range("B145").Value = 7532414
$ws.Range("E145").Value = 'Independiente Petrolero'
$ws.Range("F145").Value = 'Real Santa Cruz'
$ws.Range("G145").Value = 1
$ws.Range("H145").Value = 0
$ws.Range("I145").Value = 1
$ws.Range("K145").Value = 'H'
$ws.Range("L145").Value = 1.571
$ws.Range("M145").Value = 3.75
$ws.Range("N145").Value = 5
$ws.Range("O145").Value = 1.3
$ws.Range("P145").Value = 5
$ws.Range("Q145").Value = 11
$ws.Range("R145").Value = -1.75
$ws.Range("S145").Value = 2
$ws.Range("T145").Value = 1.8
$ws.Range("U145").Value = 3
$ws.Range("V145").Value = 1.85
$ws.Range("W145").Value = 1.95
$ws.Range("X145").Value = 0.3
$ws.Range("Z145").Value = -1
$ws.Range("AB145").Value = 0.8
$ws.Range("AD145").Value = 0.95

# Row 148: swap in match data that was previously on row 150 (id in column A is unchanged)
$ws.Range("B148").Value = 7532421
$ws.Range("E148").Value = 'Guabira'
$ws.Range("F148").Value = 'Independiente Petrolero'
$ws.Range("G148").Value = 2
$ws.Range("I148").Value = 1
$ws.Range("L148").Value = 1.4
$ws.Range("M148").Value = 4.5
$ws.Range("N148").Value = 7.5
$ws.Range("O148").Value = 1.333
$ws.Range("P148").Value = 5.5
$ws.Range("Q148").Value = 9.5
$ws.Range("S148").Value = 1.85
$ws.Range("T148").Value = 1.95
$ws.Range("U148").Value = 3
$ws.Range("V148").Value = 1.825
$ws.Range("W148").Value = 1.975
$ws.Range("X148").Value = 0.333
$ws.Range("AA148").Value = 0.8500000000000001
$ws.Range("AC148").Value = -1
$ws.Range("AD148").Value = 0.9750000000000001

# Row 149: swap in match data that was previously on row 148 (id in column A is unchanged)
$ws.Range("B149").Value = 7532420
$ws.Range("E149").Value = 'Club Aurora'
$ws.Range("F149").Value = 'Vaca Diez'
$ws.Range("L149").Value = 1.333
$ws.Range("M149").Value = 5
$ws.Range("N149").Value = 8
$ws.Range("O149").Value = 1.3
$ws.Range("P149").Value = 6.5
$ws.Range("Q149").Value = 7
$ws.Range("R149").Value = -1.5
$ws.Range("S149").Value = 1.8
$ws.Range("T149").Value = 2
$ws.Range("U149").Value = 3.25
$ws.Range("X149").Value = 0.3
$ws.Range("AA149").Value = 0.8
$ws.Range("AC149").Value = -0.5
$ws.Range("AD149").Value = 0.425

# Row 150: swap in match data that was previously on row 149 (id in column A is unchanged)
$ws.Range("B150").Value = 7532419
$ws.Range("E150").Value = 'Oriente Petrolero'
$ws.Range("F150").Value = 'Jorge Wilstermann'
$ws.Range("G150").Value = 3
$ws.Range("I150").Value = 2
$ws.Range("L150").Value = 2.2
$ws.Range("M150").Value = 2.5
$ws.Range("N150").Value = 4.5
$ws.Range("O150").Value = 2.375
$ws.Range("P150").Value = 2.45
$ws.Range("Q150").Value = 4.5
$ws.Range("R150").Value = -0.25
$ws.Range("S150").Value = 1.9
$ws.Range("T150").Value = 1.9
$ws.Range("U150").Value = 2
$ws.Range("V150").Value = 1.95
$ws.Range("W150").Value = 1.85
$ws.Range("X150").Value = 1.375
$ws.Range("AA150").Value = 0.8999999999999999
$ws.Range("AC150").Value = 0.95
$ws.Range("AD150").Value = -1
